$wb = $excel.ActiveWorkbook

# New values for column B (dates switched from dot separators to dash separators)
$dates = @("2024-01-27", "2024-01-27", "2024-02-01", "2024-02-15", "2024-03-09", "2024-03-16")

# New values for column C (names shifted up one row vs. the old layout; row 7 ends up blank)
$names = @("南宁·第一届异次元动漫嘉年华", "南宁·第五届小蜜蜂动漫嘉年华", "南宁·AP动漫游戏嘉年华", "南宁·草莓动漫节", "南宁·2024良牙动漫冬季盛典（冬典）", "")

# New values for column F (attendance counts bumped up)
$counts = @(2048, 605, 1389, 6849, 171, 84)

# Apply the same edits to both sheets that contain this data table:
#   sheet index 1 = "展览", sheet index 4 = "全部类型"
$sheetIndexes = @(1, 4)

foreach ($si in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($si)

    for ($i = 0; $i -lt 6; $i++) {
        $row = $i + 2

        # Column B: keep as literal text, not an auto-converted date serial
        $bCell = $ws.Range("B$row")
        $bCell.NumberFormat = "@"
        $bCell.Value = $dates[$i]

        # Column C: new event name (row 7 becomes empty)
        $ws.Range("C$row").Value = $names[$i]

        # Column F: updated numeric count
        $ws.Range("F$row").Value = $counts[$i]

        # Column I: link removed/cleared
        $ws.Range("I$row").Value = ""
    }
}
